# Rotate the "Recorded By" (column G) names left by one position for any
# cell whose value contains multiple comma-separated names, e.g.
#   "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
# Cells with a single name (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2

    if ($null -ne $val -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
            $cell.Value2 = [string]::Join(", ", $rotated)
        }
    }
}
